$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.03641152381897
$ws.Range("B1").Value = 2.049863576889038
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.844972968101501
$ws.Range("E1").Value = 1.179885029792786
